$p = $ppt.ActivePresentation

# "MEM score" slide (current slide 2) is duplicated; PowerPoint inserts the
# copy immediately after the source slide, i.e. before the "Method A..."
# slide - matching the new slide (id 259) that lands between id 258 and
# id 257 in the target deck.
$srcSlide = $p.Slides.Item(2)
$dupRange = $srcSlide.Duplicate()
$newSlide = $dupRange.Item(1)

# Use the existing "TextBox 8" (the URL textbox) on the new slide as a
# template so the new shape inherits the same body/line formatting
# (wrap="square", rtlCol, spAutoFit, noFill, lstStyle, etc.), then
# reposition/rename/retext it into the new explanatory note.
$templateShapeRange = $newSlide.Shapes.Item(4).Duplicate()
$noteShape = $templateShapeRange.Item(1)

$noteShape.Name = "TextBox 4"

$noteShape.TextFrame.TextRange.Text = "The + or " + [char]8722 + " value provided along with the marker name is converted to a " + [char]8722 + "10 to +10 scale and rounded to the nearest integer. As implemented here, the maximum of the scale was set using the highest absolute value MEM score observed across all markers and populations. All values in the matrix are divided by this maximum value and multiplied by 10 to achieve the " + [char]8722 + "10 to +10 scaling. "

$noteShape.Left = 1045029 / 12700.0
$noteShape.Top = 1883228 / 12700.0
$noteShape.Width = 7935685 / 12700.0
$noteShape.Height = 1477328 / 12700.0
